# Auto-generated: apply scheduled market-data refresh to Leve profit sheets.
# For each touched cell we either set the new literal value (changed/added cells)
# or clear the cell entirely (cells removed by the refresh, e.g. when a profit
# column no longer applies for a leve).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 282.7143
$ws.Range("J19").Value = 300.5
$ws.Range("L19").Value = 300.5
$ws.Range("N19").Value = -650.5
$ws.Range("H40").Value = 7692.154
$ws.Range("I40").Value = 3999
$ws.Range("K40").Value = 3999
$ws.Range("M40").Value = -3824
$ws.Range("H64").Value = 11946.789
$ws.Range("I64").Value = 7582.9165
$ws.Range("K64").Value = 7582.9165
$ws.Range("M64").Value = -7334.9165
$ws.Range("H67").Value = 11946.789
$ws.Range("I67").Value = 7582.9165
$ws.Range("K67").Value = 7582.9165
$ws.Range("M67").Value = -6724.9165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1057
$ws.Range("I2").Value = 1057
$ws.Range("K2").Value = 1057
$ws.Range("M2").Value = -944
$ws.Range("H61").Value = 4344.5454
$ws.Range("I61").Value = 1798.6666
$ws.Range("K61").Value = 1798.6666
$ws.Range("M61").Value = -1586.6666
$ws.Range("H95").Value = 24950
$ws.Range("J95").Value = 24950
$ws.Range("L95").Value = 24950
$ws.Range("N95").Value = -30442
$ws.Range("H116").Value = 1057
$ws.Range("I116").Value = 1057
$ws.Range("K116").Value = 1057
$ws.Range("M116").Value = 1237
$ws.Range("H122").Value = 1187.5
$ws.Range("I122").Value = 1187.5
$ws.Range("K122").Value = 3562.5
$ws.Range("M122").Value = -1112.5
$ws.Range("H129").Value = 60000
$ws.Range("J129").Value = 60000
$ws.Range("L129").Value = 60000
$ws.Range("N129").Value = -70000
$ws.Range("H132").Value = 977.6316
$ws.Range("I132").Value = 977.6316
$ws.Range("K132").Value = 2932.8948
$ws.Range("M132").Value = -402.8948
$ws.Range("H136").Value = 4344.5454
$ws.Range("I136").Value = 1798.6666
$ws.Range("K136").Value = 5395.9998
$ws.Range("M136").Value = -2845.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1057
$ws.Range("I3").Value = 1057
$ws.Range("K3").Value = 1057
$ws.Range("M3").Value = -943
$ws.Range("H134").Value = 1537.7333
$ws.Range("I134").Value = 1005.1539
$ws.Range("J134").Value = 4999.5
$ws.Range("K134").Value = 3015.4617
$ws.Range("L134").Value = 14998.5
$ws.Range("M134").Value = -480.4616999999998
$ws.Range("N134").Value = -20068.5
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 849.75
$ws.Range("I58").Value = 849.75
$ws.Range("K58").Value = 849.75
$ws.Range("M58").Value = -646.75
$ws.Range("H132").Value = 8692.412
$ws.Range("I132").Value = 1538.8572
$ws.Range("J132").Value = 13699.9
$ws.Range("K132").Value = 4616.571599999999
$ws.Range("L132").Value = 41099.7
$ws.Range("M132").Value = -2086.571599999999
$ws.Range("N132").Value = -46159.7
$ws.Range("H134").Value = 2435.5144
$ws.Range("I134").Value = 2262.7742
$ws.Range("K134").Value = 6788.3226
$ws.Range("M134").Value = -4253.3226
$ws.Range("H136").Value = 849.75
$ws.Range("I136").Value = 849.75
$ws.Range("K136").Value = 2549.25
$ws.Range("M136").Value = 0.75
$ws.Range("H137").Value = 70780
$ws.Range("J137").Value = 70780
$ws.Range("L137").Value = 70780
$ws.Range("N137").Value = -80980

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4625.25
$ws.Range("I3").Value = 4625.25
$ws.Range("K3").Value = 13875.75
$ws.Range("M3").Value = -13763.75
$ws.Range("H49").Value = 2650
$ws.Range("I49").Value = 2650
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 7950
$ws.Range("L49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -7794
$ws.Range("H50").Value = 379.2857
$ws.Range("I50").Value = 379.2857
$ws.Range("K50").Value = 1137.8571
$ws.Range("M50").Value = -656.8571000000002
$ws.Range("H53").Value = 379.2857
$ws.Range("I53").Value = 379.2857
$ws.Range("K53").Value = 1137.8571
$ws.Range("M53").Value = -656.8571000000002
$ws.Range("H133").Value = 4000.5
$ws.Range("I133").Value = 4000.5
$ws.Range("K133").Value = 12001.5
$ws.Range("M133").Value = -6941.5
$ws.Range("H136").Value = 8666.666999999999
$ws.Range("I136").Value = 8666.666999999999
$ws.Range("K136").Value = 26000.001
$ws.Range("M136").Value = -20900.001
$ws.Range("H140").Value = 2699.5
$ws.Range("I140").Value = 899
$ws.Range("J140").Value = 4500
$ws.Range("K140").Value = 2697
$ws.Range("L140").Value = 13500
$ws.Range("M140").Value = 2483
$ws.Range("N140").Value = -23860
$ws.Range("H141").Value = 1833.3334
$ws.Range("I141").Value = 1833.3334
$ws.Range("K141").Value = 5500.0002
$ws.Range("M141").Value = -320.0002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 19000
$ws.Range("J33").Value = 19000
$ws.Range("L33").Value = 19000
$ws.Range("N33").Value = -19504
$ws.Range("H63").Value = 20000
$ws.Range("J63").Value = 20000
$ws.Range("L63").Value = 20000
$ws.Range("N63").Value = -21372
$ws.Range("H66").Value = 20000
$ws.Range("J66").Value = 20000
$ws.Range("L66").Value = 60000
$ws.Range("N66").Value = -66864
$ws.Range("H122").Value = 5225.4546
$ws.Range("I122").Value = 5499.6665
$ws.Range("J122").Value = 3991.5
$ws.Range("K122").Value = 16498.9995
$ws.Range("L122").Value = 11974.5
$ws.Range("M122").Value = -14048.9995
$ws.Range("N122").Value = -16874.5
$ws.Range("H126").Value = 374486600
$ws.Range("I126").Value = 374486600
$ws.Range("K126").Value = 1123459800
$ws.Range("M126").Value = -1123457330
$ws.Range("H132").Value = 3781.2144
$ws.Range("J132").Value = 5438.3335
$ws.Range("L132").Value = 16315.0005
$ws.Range("N132").Value = -21375.0005
$ws.Range("H134").Value = 74993
$ws.Range("J134").Value = 74993
$ws.Range("L134").Value = 224979
$ws.Range("N134").Value = -230049

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4261
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4261
$ws.Range("K22").Value = 0
$ws.Range("L22").ClearContents()
$ws.Range("M22").Value = 4261
$ws.Range("N22").Value = -4851
$ws.Range("H27").Value = 4261
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4261
$ws.Range("K27").Value = 0
$ws.Range("L27").ClearContents()
$ws.Range("M27").Value = 4261
$ws.Range("N27").Value = -4475
$ws.Range("H40").Value = 11320.695
$ws.Range("I40").Value = 11434
$ws.Range("J40").Value = 10999.667
$ws.Range("K40").Value = 11434
$ws.Range("L40").Value = 10999.667
$ws.Range("M40").Value = -11298
$ws.Range("N40").Value = -11271.667
$ws.Range("H122").Value = 7579.1665
$ws.Range("I122").Value = 6996.75
$ws.Range("K122").Value = 20990.25
$ws.Range("M122").Value = -18540.25
$ws.Range("H136").Value = 4187.5
$ws.Range("I136").Value = 4500
$ws.Range("J136").Value = 3875
$ws.Range("K136").Value = 13500
$ws.Range("L136").Value = 11625
$ws.Range("M136").Value = -10950
$ws.Range("N136").Value = -16725

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3329
$ws.Range("I96").Value = 3245
$ws.Range("J96").Value = 3749
$ws.Range("K96").Value = 3245
$ws.Range("L96").Value = 3749
$ws.Range("M96").Value = -1872
$ws.Range("N96").Value = -6495
$ws.Range("H132").Value = 2664.375
$ws.Range("I132").Value = 2452.0454
$ws.Range("K132").Value = 7356.1362
$ws.Range("M132").Value = -4826.1362
